$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 1917.4546
$ws.Cells.Item(33, 9).Value = 105
$ws.Cells.Item(33, 11).Value = 105
$ws.Cells.Item(33, 13).Value = 124

$ws.Cells.Item(55, 8).Value = 1350
$ws.Cells.Item(55, 10).Value = 1083.3334
$ws.Cells.Item(55, 12).Value = 1083.3334
$ws.Cells.Item(55, 14).Value = -1511.3334

$ws.Cells.Item(62, 8).Value = 2001.6666
$ws.Cells.Item(62, 9).Value = 2001.6666
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 2001.6666
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 13).Value = -1377.6666
$ws.Cells.Item(62, 14).Value = $null

$ws.Cells.Item(64, 8).Value = 94125.17999999999
$ws.Cells.Item(64, 10).Value = 3597.125
$ws.Cells.Item(64, 12).Value = 3597.125
$ws.Cells.Item(64, 14).Value = -4093.125

$ws.Cells.Item(65, 8).Value = 2001.6666
$ws.Cells.Item(65, 9).Value = 2001.6666
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 10008.333
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = -6888.333000000001
$ws.Cells.Item(65, 14).Value = $null

$ws.Cells.Item(67, 8).Value = 94125.17999999999
$ws.Cells.Item(67, 10).Value = 3597.125
$ws.Cells.Item(67, 12).Value = 3597.125
$ws.Cells.Item(67, 14).Value = -5313.125

$ws.Cells.Item(86, 8).Value = 4519.241
$ws.Cells.Item(86, 9).Value = 3880.3125
$ws.Cells.Item(86, 11).Value = 3880.3125
$ws.Cells.Item(86, 13).Value = -2757.3125

$ws.Cells.Item(89, 8).Value = 4519.241
$ws.Cells.Item(89, 9).Value = 3880.3125
$ws.Cells.Item(89, 11).Value = 19401.5625
$ws.Cells.Item(89, 13).Value = -13785.5625

$ws.Cells.Item(106, 8).Value = 2702.4285
$ws.Cells.Item(106, 9).Value = 2702.4285
$ws.Cells.Item(106, 11).Value = 2702.4285
$ws.Cells.Item(106, 13).Value = -2071.4285

$ws.Cells.Item(132, 8).Value = 4633835.5
$ws.Cells.Item(132, 9).Value = 4811974.5
$ws.Cells.Item(132, 10).Value = 2226
$ws.Cells.Item(132, 11).Value = 14435923.5
$ws.Cells.Item(132, 12).Value = 6678
$ws.Cells.Item(132, 13).Value = -14433393.5
$ws.Cells.Item(132, 14).Value = -11738

$ws.Cells.Item(137, 8).Value = 1588.9574
$ws.Cells.Item(137, 9).Value = 1199.2858
$ws.Cells.Item(137, 11).Value = 3597.8574
$ws.Cells.Item(137, 13).Value = -1047.8574

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).Value = $null

$ws.Cells.Item(32, 8).Value = 33361.906
$ws.Cells.Item(32, 9).Value = 5864.3096
$ws.Cells.Item(32, 10).Value = 129603.5
$ws.Cells.Item(32, 11).Value = 5864.3096
$ws.Cells.Item(32, 12).Value = 129603.5
$ws.Cells.Item(32, 13).Value = -5577.3096
$ws.Cells.Item(32, 14).Value = -130177.5

$ws.Cells.Item(33, 8).Value = 5000
$ws.Cells.Item(33, 9).Value = 5000
$ws.Cells.Item(33, 11).Value = 5000
$ws.Cells.Item(33, 13).Value = -4671

$ws.Cells.Item(45, 8).Value = 1569.4286
$ws.Cells.Item(45, 9).Value = 1288.1666
$ws.Cells.Item(45, 11).Value = 1288.1666
$ws.Cells.Item(45, 13).Value = -911.1666

$ws.Cells.Item(61, 8).Value = 2276.2104
$ws.Cells.Item(61, 9).Value = 1558.8235
$ws.Cells.Item(61, 11).Value = 1558.8235
$ws.Cells.Item(61, 13).Value = -1346.8235

$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).Value = $null

$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).Value = $null

$ws.Cells.Item(106, 8).Value = 39370
$ws.Cells.Item(106, 10).Value = 39370
$ws.Cells.Item(106, 12).Value = 39370
$ws.Cells.Item(106, 14).Value = -41894

$ws.Cells.Item(110, 8).Value = 45547870
$ws.Cells.Item(110, 9).Value = 83501130
$ws.Cells.Item(110, 11).Value = 83501130
$ws.Cells.Item(110, 13).Value = -83499085

$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).Value = $null

$ws.Cells.Item(132, 8).Value = 3458.1738
$ws.Cells.Item(132, 9).Value = 3855.6177
$ws.Cells.Item(132, 11).Value = 11566.8531
$ws.Cells.Item(132, 13).Value = -9036.8531

$ws.Cells.Item(136, 8).Value = 2276.2104
$ws.Cells.Item(136, 9).Value = 1558.8235
$ws.Cells.Item(136, 11).Value = 4676.470499999999
$ws.Cells.Item(136, 13).Value = -2126.470499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 64535.375
$ws.Cells.Item(20, 9).Value = 73056.14
$ws.Cells.Item(20, 11).Value = 73056.14
$ws.Cells.Item(20, 13).Value = -72809.14

$ws.Cells.Item(112, 8).Value = 35469
$ws.Cells.Item(112, 10).Value = 35469
$ws.Cells.Item(112, 12).Value = 35469
$ws.Cells.Item(112, 14).Value = -38423

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 10896.75
$ws.Cells.Item(58, 9).Value = 1559.6
$ws.Cells.Item(58, 10).Value = 21670.385
$ws.Cells.Item(58, 11).Value = 1559.6
$ws.Cells.Item(58, 12).Value = 21670.385
$ws.Cells.Item(58, 13).Value = -1356.6
$ws.Cells.Item(58, 14).Value = -22076.385

$ws.Cells.Item(136, 8).Value = 10896.75
$ws.Cells.Item(136, 9).Value = 1559.6
$ws.Cells.Item(136, 10).Value = 21670.385
$ws.Cells.Item(136, 11).Value = 4678.799999999999
$ws.Cells.Item(136, 12).Value = 65011.155
$ws.Cells.Item(136, 13).Value = -2128.799999999999
$ws.Cells.Item(136, 14).Value = -70111.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 35.842106
$ws.Cells.Item(12, 10).Value = 37.64706
$ws.Cells.Item(12, 12).Value = 112.94118
$ws.Cells.Item(12, 14).Value = -458.94118

$ws.Cells.Item(37, 8).Value = 648318.1
$ws.Cells.Item(37, 10).Value = 648318.1
$ws.Cells.Item(37, 12).Value = 1944954.3
$ws.Cells.Item(37, 14).Value = -1945178.3

$ws.Cells.Item(131, 8).Value = 814.0217
$ws.Cells.Item(131, 9).Value = 547.5
$ws.Cells.Item(131, 10).Value = 832.6163
$ws.Cells.Item(131, 11).Value = 1642.5
$ws.Cells.Item(131, 12).Value = 2497.8489
$ws.Cells.Item(131, 13).Value = 3397.5
$ws.Cells.Item(131, 14).Value = -12577.8489

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 47000
$ws.Cells.Item(64, 10).Value = 47000
$ws.Cells.Item(64, 12).Value = 47000
$ws.Cells.Item(64, 14).Value = -47496

$ws.Cells.Item(67, 8).Value = 47000
$ws.Cells.Item(67, 10).Value = 47000
$ws.Cells.Item(67, 12).Value = 47000
$ws.Cells.Item(67, 14).Value = -48716

$ws.Cells.Item(80, 8).Value = 2442.6924
$ws.Cells.Item(80, 9).Value = 2633.2222
$ws.Cells.Item(80, 10).Value = 2014
$ws.Cells.Item(80, 11).Value = 2633.2222
$ws.Cells.Item(80, 12).Value = 2014
$ws.Cells.Item(80, 13).Value = -1635.2222
$ws.Cells.Item(80, 14).Value = -4010

$ws.Cells.Item(83, 8).Value = 2442.6924
$ws.Cells.Item(83, 9).Value = 2633.2222
$ws.Cells.Item(83, 10).Value = 2014
$ws.Cells.Item(83, 11).Value = 13166.111
$ws.Cells.Item(83, 12).Value = 10070
$ws.Cells.Item(83, 13).Value = -8174.111000000001
$ws.Cells.Item(83, 14).Value = -20054

$ws.Cells.Item(103, 8).Value = 22252
$ws.Cells.Item(103, 10).Value = 22252
$ws.Cells.Item(103, 12).Value = 22252
$ws.Cells.Item(103, 14).Value = -24596

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1635.5294
$ws.Cells.Item(7, 9).Value = 1376.9333
$ws.Cells.Item(7, 11).Value = 1376.9333
$ws.Cells.Item(7, 13).Value = -1264.9333

$ws.Cells.Item(68, 8).Value = 4067.8572
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 4067.8572
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 4067.8572
$ws.Cells.Item(68, 13).Value = $null
$ws.Cells.Item(68, 14).Value = -5565.8572

$ws.Cells.Item(71, 8).Value = 4067.8572
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 4067.8572
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 20339.286
$ws.Cells.Item(71, 13).Value = $null
$ws.Cells.Item(71, 14).Value = -27827.286

$ws.Cells.Item(119, 8).Value = 39000
$ws.Cells.Item(119, 10).Value = 39000
$ws.Cells.Item(119, 12).Value = 39000
$ws.Cells.Item(119, 14).Value = -48676

$ws.Cells.Item(126, 8).Value = 1635.5294
$ws.Cells.Item(126, 9).Value = 1376.9333
$ws.Cells.Item(126, 11).Value = 4130.7999
$ws.Cells.Item(126, 13).Value = -1660.7999

$ws.Cells.Item(132, 8).Value = 4981.5356
$ws.Cells.Item(132, 9).Value = 6017.8237
$ws.Cells.Item(132, 10).Value = 3380
$ws.Cells.Item(132, 11).Value = 18053.4711
$ws.Cells.Item(132, 12).Value = 10140
$ws.Cells.Item(132, 13).Value = -15523.4711
$ws.Cells.Item(132, 14).Value = -15200

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(56, 8).Value = 31865.334
$ws.Cells.Item(56, 10).Value = 31865.334
$ws.Cells.Item(56, 12).Value = 31865.334
$ws.Cells.Item(56, 14).Value = -33293.334

$ws.Cells.Item(62, 8).Value = 2073031.4
$ws.Cells.Item(62, 9).Value = 3108282.8
$ws.Cells.Item(62, 10).Value = 2528.5715
$ws.Cells.Item(62, 11).Value = 3108282.8
$ws.Cells.Item(62, 12).Value = 2528.5715
$ws.Cells.Item(62, 13).Value = -3107658.8
$ws.Cells.Item(62, 14).Value = -3776.5715

$ws.Cells.Item(65, 8).Value = 2073031.4
$ws.Cells.Item(65, 9).Value = 3108282.8
$ws.Cells.Item(65, 10).Value = 2528.5715
$ws.Cells.Item(65, 11).Value = 15541414
$ws.Cells.Item(65, 12).Value = 12642.8575
$ws.Cells.Item(65, 13).Value = -15538294
$ws.Cells.Item(65, 14).Value = -18882.8575

$ws.Cells.Item(102, 8).Value = 39765.332
$ws.Cells.Item(102, 10).Value = 39765.332
$ws.Cells.Item(102, 12).Value = 39765.332
$ws.Cells.Item(102, 14).Value = -46255.332

$ws.Cells.Item(112, 8).Value = 21390
$ws.Cells.Item(112, 10).Value = 21390
$ws.Cells.Item(112, 12).Value = 21390
$ws.Cells.Item(112, 14).Value = -24344

$ws.Cells.Item(132, 8).Value = 2792.36
$ws.Cells.Item(132, 9).Value = 3120.7
$ws.Cells.Item(132, 11).Value = 9362.099999999999
$ws.Cells.Item(132, 13).Value = -6832.099999999999
